$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Fill in the "JQuery" topic label for the two rows that already hold
#     code/definition text in column C (rows 23 & 24), matching the pattern
#     used by the preceding JQuery rows (21-22). ---
$ws.Range("A23").Value = "JQuery"
$ws.Range("A24").Value = "JQuery"

# --- Append a brand-new row (25) describing jQuery JSONP / callback usage,
#     together with a full HTML code sample. ---
$ws.Range("A25").Value = "JQuery"

$html = @'
<!DOCTYPE html>
<html lang="zh-cn">
 <head>
  <meta charset="UTF-8">
  <meta http-equiv="X-UA-Compatible" content="IE=edge">
  <meta name="viewport" content="width=device-width, initial-scale=1">
  <meta name="description" content="Home page of Handson, a social entreprise devoted on connecting technology and social issue">
  <meta name="author" content="Myles Ieong">
  <meta name="robots" content="index, follow">
  <!-- JQuery scripts and Bootstrap Core scripts -->
  <script src="https://ajax.googleapis.com/ajax/libs/jquery/1.12.0/jquery.min.js"></script>
  <script>
  function ws_results(obj) {
   alert(obj.ResultSet.totalResultsAvailable);
  }
  $(document).ready(function(){
   $("p").click(function(){
    $(this).hide();
    jQuery.getJSON("http://api.flickr.com/services/feeds/photos_public.gne?tags=cat&tagmode=any&format=json&jsoncallback=?", 
     function(data) {
      alert("response: " + data.title);
     }
    );
   });
  });
  </script>
 </head>
 <body>
  <p>If you click on me, I will disappear.</p>
  <p>Click me away!</p>
  <p>Click me too!</p>
 </body>
</html>
'@
$html = $html -replace "`r`n", "`n"
$html = $html -replace "`n", "`r`n"

$ws.Range("C25").Value = $html
$ws.Range("B25").Value = "callback + jsonp + 跨域"

# Match the styling already used elsewhere: column C cells holding long code
# blocks wrap text; column B keeps the plain (non-wrapping) default style.
$ws.Range("C25").WrapText = $true

# Row height, consistent with every other data row on the sheet.
$ws.Rows.Item(25).RowHeight = 31.5

# --- Update the view so the new row is visible and selected, mirroring the
#     author's saved workbook state. ---
$ws.Range("B26").Select() | Out-Null
